$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new "Mercedes Rental" travel expense row.
$ws.Range("A6").Value = "Mercedes Rental"
$ws.Range("B6").Value = 2500
$ws.Range("B6").NumberFormat = $ws.Range("B5").NumberFormat

# Update the selection to match the committed view state.
$ws.Range("H9").Select()
